# Append new product row 28 to the worksheet (id 610 / SKU-610-556)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "'610"
$ws.Range("A28").ClearFormats()
$ws.Range("B28").Value = "SKU-610-556"
$ws.Range("C28").Value = "new"
$ws.Range("D28").Value = "new"
$ws.Range("F28").Value = "new"
$ws.Range("G28").Value = "new"
$ws.Range("H28").Value = "new"
$ws.Range("L28").Value = "new"
$ws.Range("M28").Value = "new"
$ws.Range("N28").Value = "new"
$ws.Range("Q28").Value = 23
$ws.Range("S28").Value = 23
$ws.Range("T28").Value = "https://res.cloudinary.com/dc3hqcovg/image/upload/v1743509747/bk6ju4wqhoysv2yy19bd.png"
$ws.Range("W28").Value = "'false"
$ws.Range("W28").ClearFormats()
$ws.Range("X28").Value = "'true"
$ws.Range("X28").ClearFormats()
$ws.Range("Y28").Value = "'false"
$ws.Range("Y28").ClearFormats()
$ws.Range("Z28").Value = "'false"
$ws.Range("Z28").ClearFormats()
$ws.Range("AB28").Value = "11ml"
$ws.Range("AC28").Value = "11ML"
$ws.Range("AG28").Value = 23
$ws.Range("AI28").Value = "new"
$ws.Range("AK28").Value = "taxable"
$ws.Range("AL28").Value = "no"
$ws.Range("AM28").Value = "no"
$ws.Range("AN28").Value = "instock"
$ws.Range("AP28").Value = "publish"
$ws.Range("AQ28").Value = "open"
$ws.Range("AS28").Value = "'false"
$ws.Range("AS28").ClearFormats()
$ws.Range("AT28").Value = "https://res.cloudinary.com/dc3hqcovg/image/upload/v1743509747/bk6ju4wqhoysv2yy19bd.png"
$ws.Range("AU28").Value = "2025-04-01T12:15:47.737Z"
$ws.Range("AV28").Value = "new"
$ws.Range("AW28").Value = 23
$ws.Range("AX28").Value = "'610"
$ws.Range("AX28").ClearFormats()
$ws.Range("AY28").Value = "2025-04-01T12:15:47.737Z"
